# Financial Values workbook edit
#
# Adds a new leading "Type" column (filled with "Search" for every rule
# row), and rewrites the Operation tag values so they are wrapped in
# "~" (gt -> ~gt~, lt -> ~lt~) plus fixes the "debtoequity" typo to
# "debttoequity". Also swaps the Beta rule's lower/upper bound values so
# they read consistently with the other "between" rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data one column to the right to make room for the
# new leading "Type" column.
$ws.Columns("A:A").Insert()

# --- New "Type" column ------------------------------------------------
# Write order matters for the shared-string table, so "Search" is
# introduced before "Type".
$ws.Range("A2").Value = "Search"
$ws.Range("A1").Value = "Type"
$ws.Range("A3").Value = "Search"
$ws.Range("A4").Value = "Search"
$ws.Range("A5").Value = "Search"
$ws.Range("A6").Value = "Search"
$ws.Range("A8").Value = "Search"

# --- Header bound labels swapped to read "lower, upper" left-to-right
$ws.Range("G1").Value = "lower bound"
$ws.Range("H1").Value = "upper bound"

# --- Operation tag fixes (now column E) + debtoequity typo fix -------
$ws.Range("E3").Value = "~lt~"
$ws.Range("E5").Value = "~gt~"
$ws.Range("E6").Value = "~gt~"
$ws.Range("E8").Value = "~gt~"
$ws.Range("D3").Value = "debttoequity"

# --- Beta "between" bounds swapped (row 4, now columns G/H) ----------
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 1.1000000000000001

# --- New column widths for the newly appended value columns ----------
$ws.Columns("G:G").ColumnWidth = 11.6
$ws.Columns("H:H").ColumnWidth = 11.45

# --- Selection, matching the authored state ---------------------------
$ws.Range("H4").Select()
